$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Rows 1-12 (1-indexed): update single summary values
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "1421"
$t.Cell(5, 1).Range.Text  = "0.00002"
$t.Cell(6, 1).Range.Text  = "0.00250"
$t.Cell(7, 1).Range.Text  = "0.00014"
$t.Cell(8, 1).Range.Text  = "0.00006"
$t.Cell(9, 1).Range.Text  = "0.00023"
$t.Cell(10, 1).Range.Text = "0.00027"
$t.Cell(11, 1).Range.Text = "0.00033"
$t.Cell(12, 1).Range.Text = "0.23372"

# Rows 44-46 (1-indexed): collapse the tab-separated per-iteration rows
# down to a single value each (the old leading number for rows
# 44/45/46 that used to live in rows 1/2/3 before they were overwritten)
$t.Cell(44, 1).Range.Text = "99.97"
$t.Cell(45, 1).Range.Text = "0.23"
$t.Cell(46, 1).Range.Text = "773"
